# remove sy-datum/sy-uzeit from demo programs
#
# Sheet2 (the active sheet) had two "system field" style demo cells:
#   B3 = sy-datum (a date serial, formatted as a date)
#   C3 = sy-uzeit (a time fraction, formatted as a time)
# These are replaced by a single literal string "Hello world" in B3,
# and column C (which only existed to host C3) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Drop the old date/time formatting on B3 before writing the new literal
# text so the cell reverts to the workbook's default (unformatted) style
# instead of carrying over the date number format.
$ws.Range("B3").ClearFormats()
$ws.Range("B3").Value = "Hello world"

# C3 (sy-uzeit) is no longer needed; remove the whole column so the
# sheet's dimension/col list shrinks back down to just column B.
$ws.Columns("C").Delete()

# Cosmetic: the embedded header picture (on Sheet1) was re-saved under a
# new generated identifier.
$ws1 = $wb.Worksheets.Item(1)
$shp = $ws1.Shapes.Item(1)
$shp.Name = "0800271CCEE91EEC9AC683DB77615147"
